$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new experiment row after row 2 (new row 3), pushing existing
# rows 3-14 down to 4-15 ---
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3,1).Value = "32, 10"
$ws.Cells.Item(3,2).Value = 0
$ws.Cells.Item(3,3).Value = 300
$ws.Cells.Item(3,4).Value = 200
$ws.Cells.Item(3,5).Value = 200
$ws.Cells.Item(3,6).Value = 100
$ws.Cells.Item(3,8).Value = 1.4095
$ws.Cells.Item(3,9).Value = 0.725

# --- Highlight the WndSize (column E) values for all data rows with the
# orange fill ---
$ws.Range("E2:E15").Interior.Color = 49407

# --- Notes column ---
$ws.Cells.Item(1,11).Value = "Notes"
$ws.Cells.Item(3,11).Value = "Signs of overfitting- Training data reach 95% accuracy around epoch 85"
$ws.Cells.Item(2,11).Value = "Orange values may be erroneous"
$ws.Cells.Item(2,11).Interior.Color = 49407

# --- View state ---
$ws.Range("K16").Select()
